# Add an "About" sheet (ahead of all existing sheets) describing the
# framework/template, matching "Add about sheet from which name is drawn".

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the first existing sheet so it becomes the
# first tab in the workbook.
$firstSheet = $wb.Worksheets.Item(1)
$about = $wb.Worksheets.Add($firstSheet)
$about.Name = "About"

# Header row
$about.Range("A1").Value = "Name"
$about.Range("B1").Value = "Description"
$about.Range("A1:B1").Font.Bold = $true

# Data row
$about.Range("A2").Value = "UDT"
$about.Range("B2").Value = "UDT example"
$about.Range("A2:B2").VerticalAlignment = -4160

# Leave the final selection on B2, matching the authored workbook, and make
# sure the About sheet is the active/selected tab.
[void]$about.Range("B2").Select()
$about.Activate()
